$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.42'
$ws.Range("G2").Value = '3'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.01'
$ws.Range("G3").Value = '3'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.337'
$ws.Range("G4").Value = '3'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05772'
$ws.Range("G5").Value = '3'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.477'
$ws.Range("G6").Value = '3'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.333'
$ws.Range("G7").Value = '3'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8108'
$ws.Range("G8").Value = '3'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9205'
$ws.Range("G9").Value = '3'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1403'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").Value = '3'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07351'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").Value = '3'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03166'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G12").Value = '3'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03062'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").Value = '3'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09345'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").Value = '3'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.872'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = '3'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001540'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("G16").Value = '3'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04764'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").Value = '3'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006042'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("G18").Value = '3'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005993'
$ws.Range("G19").Value = '3'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001279'
$ws.Range("G20").Value = '3'
$ws.Range("G21").Value = '3'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008838'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("G22").Value = '3'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.606'
$ws.Range("G23").Value = '3'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.138'
$ws.Range("G24").Value = '3'
$ws.Range("G25").Value = '3'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1318'
$ws.Range("G26").Value = '3'
$ws.Range("G27").Value = '3'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002354'
$ws.Range("G28").Value = '3'
$ws.Range("G29").Value = '3'
$ws.Range("G30").Value = '3'
$ws.Range("G31").Value = '3'
$ws.Range("G32").Value = '3'
$ws.Range("G33").Value = '3'
$ws.Range("G34").Value = '3'
$ws.Range("G35").Value = '3'
$ws.Range("G36").Value = '3'
$ws.Range("G37").Value = '3'
$ws.Range("G38").Value = '3'
$ws.Range("G39").Value = '3'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03799'
$ws.Range("G40").Value = '3'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006387'
$ws.Range("G41").Value = '3'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1053'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '3'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002762'
$ws.Range("E43").Value = '42CEJICEJIWorstin24h'
$ws.Range("G43").Value = '3'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008437'
$ws.Range("G44").Value = '3'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005341'
$ws.Range("G45").Value = '3'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000752'
$ws.Range("G46").Value = '3'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6913'
$ws.Range("G47").Value = '3'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.001852'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("G48").Value = '3'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002104'
$ws.Range("G49").Value = '3'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002004'
$ws.Range("G50").Value = '3'
$ws.Range("G51").Value = '3'
